# Fix data import: header "Deskripsi" column was replaced with "Kategori"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header in column C (was "Deskripsi", now "Kategori")
$ws.Range("C1").Value = "Kategori"

# Move the active selection to D6 (matches the saved cursor position)
$ws.Range("D6").Select()
